# Updated cryptos list on Mon Feb 12 10:58:15 UTC 2024 with GitHub Actions
# Refresh the live Price (D) / Volume(1h) (E) columns, and fix the row 42/43
# coin swap (EnergySwap <-> WEMIXToken), to match the latest coinranking.com pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "47.861.68"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.478.95"
$ws.Range("E3").Value = "  -1.49%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5: BNB
$ws.Range("D5").Value = "'316.81"
$ws.Range("E5").Value = "  -1.55%  "

# Row 6: Solana
$ws.Range("D6").Value = "'104.58"
$ws.Range("E6").Value = "  -3.91%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.517"
$ws.Range("E7").Value = "  -2.69%  "

# Row 8: USDC
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.07%  "

# Row 9: Cardano
$ws.Range("D9").Value = "'0.534"

# Row 10: Avalanche
$ws.Range("D10").Value = "'38.75"
$ws.Range("E10").Value = "  -4.02%  "

# Row 11: Chainlink
$ws.Range("D11").Value = "'20.37"
$ws.Range("E11").Value = "  -0.90%  "

# Row 12: Dogecoin
$ws.Range("D12").Value = "'0.0797"
$ws.Range("E12").Value = "  -3.07%  "

# Row 13: TRON
$ws.Range("E13").Value = "  +0.22%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "'7.01"

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.869.77"
$ws.Range("E15").Value = "  -1.43%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "2.506.84"
$ws.Range("E16").Value = "  -0.43%  "

# Row 17: Polygon
$ws.Range("D17").Value = "'0.820"
$ws.Range("E17").Value = "  -3.95%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "47.806.62"
$ws.Range("E18").Value = "  -0.14%  "

# Row 19: ImmutableX
$ws.Range("D19").Value = "'2.92"
$ws.Range("E19").Value = "  +8.11%  "

# Row 20: InternetComputer(DFINITY)
$ws.Range("D20").Value = "'12.61"
$ws.Range("E20").Value = "  -4.83%  "

# Row 21: Uniswap
$ws.Range("E21").Value = "  -1.44%  "

# Row 22: ShibaInu
$ws.Range("D22").Value = "0.0₃0924"
$ws.Range("E22").Value = "  -2.10%  "

# Row 23: BitcoinCash
$ws.Range("D23").Value = "'278.90"
$ws.Range("E23").Value = "  +5.34%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'70.71"
$ws.Range("E24").Value = "  -1.64%  "

# Row 25: PancakeSwap
$ws.Range("D25").Value = "'2.48"
$ws.Range("E25").Value = "  -3.13%  "

# Row 26: Dai
$ws.Range("E26").Value = "  +0.19%  "

# Row 27: EthereumClassic
$ws.Range("D27").Value = "'25.59"
$ws.Range("E27").Value = "  -1.73%  "

# Row 28: Toncoin
$ws.Range("D28").Value = "'2.21"
$ws.Range("E28").Value = "  +0.21%  "

# Row 29: Cosmos
$ws.Range("D29").Value = "'9.55"
$ws.Range("E29").Value = "  -5.40%  "

# Row 30: Kaspa
$ws.Range("E30").Value = "  -4.11%  "

# Row 31: InjectiveProtocol
$ws.Range("D31").Value = "'34.45"
$ws.Range("E31").Value = "  -4.00%  "

# Row 32: OKB
$ws.Range("E32").Value = "  -1.03%  "

# Row 33: FirstDigitalUSD
$ws.Range("E33").Value = "  -0.12%  "

# Row 34: Celestia
$ws.Range("D34").Value = "'18.83"
$ws.Range("E34").Value = "  -4.67%  "

# Row 35: Filecoin
$ws.Range("E35").Value = "  -2.95%  "

# Row 36: Hedera
$ws.Range("D36").Value = "'0.0765"
$ws.Range("E36").Value = "  -2.85%  "

# Row 37: ARBITRUM
$ws.Range("E37").Value = "  -2.47%  "

# Row 38: RenderToken
$ws.Range("E38").Value = "  -4.02%  "

# Row 39: LidoDAOToken
$ws.Range("E39").Value = "  -4.34%  "

# Row 40: Monero
$ws.Range("D40").Value = "'122.03"
$ws.Range("E40").Value = "  +1.99%  "

# Row 41: Stellar
$ws.Range("E41").Value = "  -1.64%  "

# Row 42: EnergySwap
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'2.20"
$ws.Range("E42").Value = "  +0.06%  "

# Row 43: WEMIXToken
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'21.87"
$ws.Range("E43").Value = "  -0.90%  "

# Row 44: VeChain
$ws.Range("D44").Value = "'0.0298"
$ws.Range("E44").Value = "  -0.40%  "

# Row 45: Maker
$ws.Range("D45").Value = "1.992.37"
$ws.Range("E45").Value = "  -0.94%  "

# Row 46: NEARProtocol
$ws.Range("D46").Value = "'3.13"
$ws.Range("E46").Value = "  -0.46%  "

# Row 47: Stacks
$ws.Range("E47").Value = "  -0.85%  "

# Row 48: ApeXProtocol
$ws.Range("E48").Value = "  -3.40%  "

# Row 49: FraxShare
$ws.Range("D49").Value = "'8.92"
$ws.Range("E49").Value = "  -1.95%  "

# Row 50: THORChain
$ws.Range("D50").Value = "'5.11"
$ws.Range("E50").Value = "  -1.63%  "

# Row 51: BitcoinSV
$ws.Range("D51").Value = "'78.89"
$ws.Range("E51").Value = "  +0.38%  "
